$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.008.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.671.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.513"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.08"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0891"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.908.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.673.24"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.75"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.77%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.038.18"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.94"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0735"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.72"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.45"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.23"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.24"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.30"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.87"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.451.96"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.14%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +5.46%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.568"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.57%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +11.65%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.75"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.816.94"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.26"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.33%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0508"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.94%  "
